$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 186.96428
$ws.Range("I33").Value = 132.18182
$ws.Range("K33").Value = 132.18182
$ws.Range("M33").Value = 96.81818000000001
$ws.Range("H41").Value = 1697.1177
$ws.Range("I41").Value = 1882.625
$ws.Range("J41").Value = 1532.2222
$ws.Range("K41").Value = 1882.625
$ws.Range("L41").Value = 1532.2222
$ws.Range("M41").Value = -1442.625
$ws.Range("N41").Value = -2412.2222
$ws.Range("H43").Value = 9259934
$ws.Range("I43").Value = 450.33334
$ws.Range("J43").Value = 18519418
$ws.Range("K43").Value = 450.33334
$ws.Range("L43").Value = 18519418
$ws.Range("M43").Value = -381.33334
$ws.Range("N43").Value = -18519556
$ws.Range("H58").Value = 3200.8635
$ws.Range("I58").Value = 484
$ws.Range("J58").Value = 3999.9412
$ws.Range("K58").Value = 1452
$ws.Range("L58").Value = 11999.8236
$ws.Range("M58").Value = -1302
$ws.Range("N58").Value = -12299.8236
$ws.Range("H80").Value = 775
$ws.Range("I80").Value = 500
$ws.Range("K80").Value = 1500
$ws.Range("M80").Value = -502
$ws.Range("H83").Value = 775
$ws.Range("I83").Value = 500
$ws.Range("K83").Value = 4500
$ws.Range("M83").Value = 492
$ws.Range("H106").Value = 2816.5833
$ws.Range("I106").Value = 2816.5833
$ws.Range("K106").Value = 2816.5833
$ws.Range("M106").Value = -2185.5833
$ws.Range("H113").Value = 66669500
$ws.Range("J113").Value = 4500
$ws.Range("L113").Value = 4500
$ws.Range("N113").Value = -11008
$ws.Range("H114").Value = 38000
$ws.Range("J114").Value = 38000
$ws.Range("L114").Value = 38000
$ws.Range("N114").Value = -46678
$ws.Range("H115").Value = 1564.1666
$ws.Range("I115").Value = 954
$ws.Range("K115").Value = 2862
$ws.Range("M115").Value = -1295
$ws.Range("H118").Value = 2670
$ws.Range("I118").Value = 1922.5
$ws.Range("K118").Value = 5767.5
$ws.Range("M118").Value = -4110.5
$ws.Range("H135").Value = 230.6
$ws.Range("I135").Value = 181.8
$ws.Range("K135").Value = 1636.2
$ws.Range("M135").Value = 898.8
$ws.Range("H137").Value = 2024.1613
$ws.Range("I137").Value = 1738.9474
$ws.Range("J137").Value = 2475.75
$ws.Range("K137").Value = 5216.8422
$ws.Range("L137").Value = 7427.25
$ws.Range("M137").Value = -2666.8422
$ws.Range("N137").Value = -12527.25

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3028.356
$ws.Range("I32").Value = 3137.7886
$ws.Range("K32").Value = 3137.7886
$ws.Range("M32").Value = -2850.7886
$ws.Range("H45").Value = 1541.091
$ws.Range("I45").Value = 1541.091
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1541.091
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = $null
$ws.Range("N45").Value = -1164.091
$ws.Range("H61").Value = 1368.8462
$ws.Range("I61").Value = 1116.4546
$ws.Range("J61").Value = 2757
$ws.Range("K61").Value = 1116.4546
$ws.Range("L61").Value = 2757
$ws.Range("M61").Value = -904.4546
$ws.Range("N61").Value = -3181
$ws.Range("H74").Value = 561.125
$ws.Range("I74").Value = 686.4545000000001
$ws.Range("J74").Value = 285.4
$ws.Range("K74").Value = 686.4545000000001
$ws.Range("L74").Value = 285.4
$ws.Range("M74").Value = 187.5454999999999
$ws.Range("N74").Value = -2033.4
$ws.Range("H77").Value = 561.125
$ws.Range("I77").Value = 686.4545000000001
$ws.Range("J77").Value = 285.4
$ws.Range("K77").Value = 3432.2725
$ws.Range("L77").Value = 1427
$ws.Range("M77").Value = 935.7275
$ws.Range("N77").Value = -10163
$ws.Range("H132").Value = 3938.5
$ws.Range("I132").Value = 3913.875
$ws.Range("K132").Value = 11741.625
$ws.Range("M132").Value = -9211.625
$ws.Range("H135").Value = 17971.5
$ws.Range("J135").Value = 17971.5
$ws.Range("L135").Value = 17971.5
$ws.Range("N135").Value = -28111.5
$ws.Range("H136").Value = 1368.8462
$ws.Range("I136").Value = 1116.4546
$ws.Range("J136").Value = 2757
$ws.Range("K136").Value = 3349.3638
$ws.Range("L136").Value = 8271
$ws.Range("M136").Value = -799.3638000000001
$ws.Range("N136").Value = -13371

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H52").Value = 29560
$ws.Range("J52").Value = 29560
$ws.Range("L52").Value = 29560
$ws.Range("N52").Value = -30086
$ws.Range("H54").Value = 4747.4
$ws.Range("I54").Value = 1559.25
$ws.Range("J54").Value = 17500
$ws.Range("K54").Value = 1559.25
$ws.Range("L54").Value = 17500
$ws.Range("M54").Value = -1075.25
$ws.Range("N54").Value = -18468
$ws.Range("H116").Value = 40556.5
$ws.Range("J116").Value = 40556.5
$ws.Range("L116").Value = 40556.5
$ws.Range("N116").Value = -49734.5
$ws.Range("H117").Value = 45000
$ws.Range("J117").Value = 45000
$ws.Range("L117").Value = 45000
$ws.Range("N117").Value = -54178
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = $null
$ws.Range("N119").Value = 0
$ws.Range("H120").Value = 30000
$ws.Range("J120").Value = 30000
$ws.Range("L120").Value = 30000
$ws.Range("N120").Value = -39676
$ws.Range("H121").Value = 29560
$ws.Range("J121").Value = 29560
$ws.Range("L121").Value = 29560
$ws.Range("N121").Value = -33054
$ws.Range("H134").Value = 8574
$ws.Range("I134").Value = 1365.7273
$ws.Range("J134").Value = 35004.332
$ws.Range("K134").Value = 4097.1819
$ws.Range("L134").Value = 105012.996
$ws.Range("M134").Value = -1562.1819
$ws.Range("N134").Value = -110082.996

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1113.6165
$ws.Range("I31").Value = 834.525
$ws.Range("J31").Value = 1451.909
$ws.Range("K31").Value = 834.525
$ws.Range("L31").Value = 1451.909
$ws.Range("M31").Value = -539.525
$ws.Range("N31").Value = -2041.909
$ws.Range("H34").Value = 1113.6165
$ws.Range("I34").Value = 834.525
$ws.Range("J34").Value = 1451.909
$ws.Range("K34").Value = 834.525
$ws.Range("L34").Value = 1451.909
$ws.Range("M34").Value = -632.525
$ws.Range("N34").Value = -1855.909
$ws.Range("H58").Value = 967.25
$ws.Range("I58").Value = 1048.2858
$ws.Range("J58").Value = 400
$ws.Range("K58").Value = 1048.2858
$ws.Range("L58").Value = 400
$ws.Range("M58").Value = -845.2858000000001
$ws.Range("N58").Value = -806
$ws.Range("H132").Value = 5769.6
$ws.Range("J132").Value = 5999
$ws.Range("L132").Value = 17997
$ws.Range("N132").Value = -23057
$ws.Range("H134").Value = 6804015
$ws.Range("I134").Value = 8334495
$ws.Range("J134").Value = 1880.8889
$ws.Range("K134").Value = 25003485
$ws.Range("L134").Value = 5642.6667
$ws.Range("M134").Value = -25000950
$ws.Range("N134").Value = -10712.6667
$ws.Range("H136").Value = 967.25
$ws.Range("I136").Value = 1048.2858
$ws.Range("J136").Value = 400
$ws.Range("K136").Value = 3144.8574
$ws.Range("L136").Value = 1200
$ws.Range("M136").Value = -594.8574000000003
$ws.Range("N136").Value = -6300

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 293.57144
$ws.Range("I50").Value = 318.33334
$ws.Range("J50").Value = 275
$ws.Range("K50").Value = 955.0000200000001
$ws.Range("L50").Value = 825
$ws.Range("M50").Value = -474.0000200000001
$ws.Range("N50").Value = -1787
$ws.Range("H53").Value = 293.57144
$ws.Range("I53").Value = 318.33334
$ws.Range("J53").Value = 275
$ws.Range("K53").Value = 955.0000200000001
$ws.Range("L53").Value = 825
$ws.Range("M53").Value = -474.0000200000001
$ws.Range("N53").Value = -1787
$ws.Range("H68").Value = 2230.0566
$ws.Range("I68").Value = 550
$ws.Range("J68").Value = 2295.9412
$ws.Range("K68").Value = 1650
$ws.Range("L68").Value = 6887.823600000001
$ws.Range("M68").Value = -839
$ws.Range("N68").Value = -8509.8236
$ws.Range("H71").Value = 2230.0566
$ws.Range("I71").Value = 550
$ws.Range("J71").Value = 2295.9412
$ws.Range("K71").Value = 4950
$ws.Range("L71").Value = 20663.4708
$ws.Range("M71").Value = -894
$ws.Range("N71").Value = -28775.4708
$ws.Range("H86").Value = 1500
$ws.Range("I86").Value = 1500
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 4500
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = $null
$ws.Range("N86").Value = -3314
$ws.Range("H89").Value = 1500
$ws.Range("I89").Value = 1500
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 13500
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = $null
$ws.Range("N89").Value = -7572
$ws.Range("H125").Value = 4954.5454
$ws.Range("J125").Value = 5500
$ws.Range("L125").Value = 16500
$ws.Range("N125").Value = -26340
$ws.Range("H137").Value = 15368.9
$ws.Range("I137").Value = 2718
$ws.Range("J137").Value = 28019.8
$ws.Range("K137").Value = 8154
$ws.Range("L137").Value = 84059.39999999999
$ws.Range("M137").Value = -3054
$ws.Range("N137").Value = -94259.39999999999
$ws.Range("H140").Value = 27198.045
$ws.Range("I140").Value = 37600.902
$ws.Range("J140").Value = 2391.2307
$ws.Range("K140").Value = 112802.706
$ws.Range("L140").Value = 7173.6921
$ws.Range("M140").Value = -107622.706
$ws.Range("N140").Value = -17533.6921

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 53800
$ws.Range("J104").Value = 53800
$ws.Range("L104").Value = 53800
$ws.Range("N104").Value = -60788
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = $null
$ws.Range("N105").Value = 0
$ws.Range("H132").Value = 3062.5217
$ws.Range("I132").Value = 2690.647
$ws.Range("J132").Value = 4116.1665
$ws.Range("K132").Value = 8071.941
$ws.Range("L132").Value = 12348.4995
$ws.Range("M132").Value = -5541.941
$ws.Range("N132").Value = -17408.4995

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 7570.8823
$ws.Range("I136").Value = 9125.923000000001
$ws.Range("J136").Value = 2517
$ws.Range("K136").Value = 27377.769
$ws.Range("L136").Value = 7551
$ws.Range("M136").Value = -24827.769
$ws.Range("N136").Value = -12651

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1296.7407
$ws.Range("I136").Value = 434.75
$ws.Range("J136").Value = 1659.6842
$ws.Range("K136").Value = 1304.25
$ws.Range("L136").Value = 4979.0526
$ws.Range("M136").Value = 1245.75
$ws.Range("N136").Value = -10079.0526
